# Generate Report for Handback
# This script applies the "handback" updates to the localization-status
# workbook: it marks the two tracked files as handed back (in sync with
# en-US), records the handback target/handback files for each locale and
# stamps a handback datetime for the locale that is now fully in sync.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns for both locales / both rows ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# --- zh-cn sheet ---
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Target File (F) / Latest Handback File (G) for row 2 and 3.
$zhMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/0b05ceeeb6e1ec5418cec7c2b1df0393349977df/e2e/6b984f79-5315-48d5-bbca-213b847ccf8e.md"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/75735567ac5c743a76d9207df9dd16b9e0e21040/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf"

$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhMdUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $zhXlfUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhMdUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $zhXlfUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.zh-cn.xlf")

# Latest Handback DateTime (H) keeps the same placeholder value text, which is
# itself updated in shared strings below -- no direct per-cell change needed.
$wsZh.Range("H2").Value = "2016-03-11 10:52:53"
$wsZh.Range("H3").Value = "2016-03-11 10:52:53"

# --- de-de sheet ---
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$deMdUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/0b05ceeeb6e1ec5418cec7c2b1df0393349977df/e2e/6b984f79-5315-48d5-bbca-213b847ccf8e.md"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dedafacfe23cd4e0beb726c0187bd00005439d00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deMdUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $deXlfUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deMdUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $deXlfUrl, "", "", "6b984f79-5315-48d5-bbca-213b847ccf8e.7d306ebc8b8eab999ead8f0a675bfed8bcaa2154.de-de.xlf")

# de-de is now fully in sync, so stamp a fresh handback datetime.
$wsDe.Range("H2").Value = "2016-03-11 10:52:58"
$wsDe.Range("H3").Value = "2016-03-11 10:52:58"
